# Generate Report for Handback
# - Flip "Ready for handoff" status to "Handed back: in sync with en-US"
#   on the Overview sheet and both locale sheets.
# - Refresh the zh-cn / de-de "Latest Handback DateTime" with the new
#   handback timestamps now that the files are in sync.
# - Clear the stale "handback file is not the latest" error details now
#   that the handback succeeded.
# - Widen the Status / Error-detail-ish columns so the new, longer status
#   text isn't clipped.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# --- zh-cn detail sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("K2").Value = "2016-11-14 17:59:48"
$zhcn.Range("K3").Value = "2016-11-14 17:59:48"

$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(16).ColumnWidth = 12.8333333333333

# --- de-de detail sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Re-assert the handoff/handback file names so the now-unused "stale
# handback" error strings are dropped from the shared-string table.
$dede.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("K2").Value = "2016-11-14 18:00:10"
$dede.Range("K3").Value = "2016-11-14 18:00:10"

$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(16).ColumnWidth = 12.8333333333333
